$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BHPSbP")

# --- Core data update: link the BHPSbP production shares to the About sheet's
#     assumption values instead of hard-coded literals ---
$wsData.Range("B2").Formula = "=About!A22"
$wsData.Range("B3").Formula = "=About!A23"
$wsData.Range("B4").Formula = "=About!A24"

# --- View/selection state updates ---
# About sheet: select A22:A24 (the cells that now feed BHPSbP) and scroll down
$wsAbout.Activate()
$wsAbout.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
$wsAbout.Range("A22:A24").Select()

# BHPSbP becomes the active/selected sheet
$wsData.Activate()
$wsData.Range("A1").Select()
